$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 294, pushing the existing
# rows 294:319 down to 296:321 (matching the rest of the diff, which is
# just those rows shifted by two positions).
$ws.Rows("294:295").Insert()

# New row 294 - Mango, "Primera" quality entry for date 44578 (2022-01-17)
$ws.Cells.Item(294, 1).Value = 3
$ws.Cells.Item(294, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(294, 3).Value = "Coquimbo"
$ws.Cells.Item(294, 4).Value = 44578
$ws.Cells.Item(294, 5).Value = 5
$ws.Cells.Item(294, 6).Value = "Fruta"
$ws.Cells.Item(294, 7).Value = 100108
$ws.Cells.Item(294, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(294, 9).Value = 100108002
$ws.Cells.Item(294, 10).Value = "Mango"
$ws.Cells.Item(294, 11).Value = "Sin especificar"
$ws.Cells.Item(294, 12).Value = "Primera"
$ws.Cells.Item(294, 13).Value = 240
$ws.Cells.Item(294, 14).Value = 6500
$ws.Cells.Item(294, 15).Value = 7000
$ws.Cells.Item(294, 16).Value = 6750
$ws.Cells.Item(294, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(294, 18).Value = "Perú"
$ws.Cells.Item(294, 19).Value = 1688
$ws.Cells.Item(294, 20).Value = 4

# New row 295 - Mango, "Segunda" quality entry for the same date 44578
$ws.Cells.Item(295, 1).Value = 3
$ws.Cells.Item(295, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(295, 3).Value = "Coquimbo"
$ws.Cells.Item(295, 4).Value = 44578
$ws.Cells.Item(295, 5).Value = 5
$ws.Cells.Item(295, 6).Value = "Fruta"
$ws.Cells.Item(295, 7).Value = 100108
$ws.Cells.Item(295, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(295, 9).Value = 100108002
$ws.Cells.Item(295, 10).Value = "Mango"
$ws.Cells.Item(295, 11).Value = "Sin especificar"
$ws.Cells.Item(295, 12).Value = "Segunda"
$ws.Cells.Item(295, 13).Value = 80
$ws.Cells.Item(295, 14).Value = 5000
$ws.Cells.Item(295, 15).Value = 5000
$ws.Cells.Item(295, 16).Value = 5000
$ws.Cells.Item(295, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(295, 18).Value = "Perú"
$ws.Cells.Item(295, 19).Value = 1250
$ws.Cells.Item(295, 20).Value = 4
